# Updated cryptos list with GitHub Actions refresh of prices / 1h volume
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column D holds price text that must stay literal text (e.g. "1.060", "10.80",
# "0.00001040") - mark it as Text before writing so Excel does not coerce it to a
# number and silently drop significant trailing/leading zeros.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '20.310.79'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '1.438.23'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.38%  '
$ws.Range('D5').Value = '0.9928'
$ws.Range('E5').Value = '  -1.16%  '
$ws.Range('D6').Value = '278.56'
$ws.Range('E6').Value = '  +1.14%  '
$ws.Range('D7').Value = '0.3703'
$ws.Range('E7').Value = '  +0.44%  '
$ws.Range('D8').Value = '0.3171'
$ws.Range('E8').Value = '  +3.55%  '
$ws.Range('D9').Value = '40.48'
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('D10').Value = '1.060'
$ws.Range('E10').Value = '  +5.65%  '
$ws.Range('D11').Value = '0.06631'
$ws.Range('E11').Value = '  +1.29%  '
$ws.Range('D12').Value = '0.9955'
$ws.Range('E12').Value = '  -0.98%  '
$ws.Range('D13').Value = '5.571'
$ws.Range('E13').Value = '  +4.11%  '
$ws.Range('D14').Value = '18.23'
$ws.Range('E14').Value = '  +6.21%  '
$ws.Range('D15').Value = '6.249'
$ws.Range('E15').Value = '  +2.05%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '0.00001040'
$ws.Range('E16').Value = '  +3.98%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '1.436.11'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').Value = '0.05761'
$ws.Range('E18').Value = '  -1.02%  '
$ws.Range('D19').Value = '0.9932'
$ws.Range('E19').Value = '  -1.07%  '
$ws.Range('D20').Value = '71.94'
$ws.Range('E20').Value = '  -5.02%  '
$ws.Range('D21').Value = '5.651'
$ws.Range('E21').Value = '  -1.13%  '
$ws.Range('D22').Value = '14.84'
$ws.Range('E22').Value = '  +3.19%  '
$ws.Range('D23').Value = '11.21'
$ws.Range('E23').Value = '  +2.37%  '
$ws.Range('D24').Value = '2.266'
$ws.Range('E24').Value = '  -2.49%  '
$ws.Range('D25').Value = '20.338.04'
$ws.Range('E25').Value = '  +0.71%  '
$ws.Range('D26').Value = '2.325'
$ws.Range('E26').Value = '  +1.78%  '
$ws.Range('D27').Value = '135.43'
$ws.Range('E27').Value = '  -4.70%  '
$ws.Range('D28').Value = '17.48'
$ws.Range('E28').Value = '  +3.02%  '
$ws.Range('D29').Value = '1.596.33'
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').Value = '112.38'
$ws.Range('E30').Value = '  +2.88%  '
$ws.Range('D31').Value = '4.004'
$ws.Range('E31').Value = '  +1.55%  '
$ws.Range('D32').Value = '5.355'
$ws.Range('E32').Value = '  -1.83%  '
$ws.Range('D33').Value = '0.8465'
$ws.Range('E33').Value = '  -8.39%  '
$ws.Range('D34').Value = '0.07815'
$ws.Range('E34').Value = '  +1.76%  '
$ws.Range('D35').Value = '1.503'
$ws.Range('E35').Value = '  +19.00%  '
$ws.Range('B36').Value = 'InternetComputer(DFINITY)'
$ws.Range('C36').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D36').Value = '4.971'
$ws.Range('E36').Value = '  +5.63%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '0.05937'
$ws.Range('E37').Value = '  +5.96%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').Value = '10.80'
$ws.Range('E38').Value = '  -2.22%  '
$ws.Range('B39').Value = 'Frax'
$ws.Range('C39').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D39').Value = '0.9938'
$ws.Range('E39').Value = '  -0.76%  '
$ws.Range('D40').Value = '7.877'
$ws.Range('E40').Value = '  -6.03%  '
$ws.Range('D41').Value = '0.02074'
$ws.Range('E41').Value = '  +3.05%  '
$ws.Range('D42').Value = '1.132'
$ws.Range('E42').Value = '  +1.09%  '
$ws.Range('D43').Value = '0.1906'
$ws.Range('E43').Value = '  -0.30%  '
$ws.Range('D44').Value = '0.5406'
$ws.Range('E44').Value = '  +2.44%  '
$ws.Range('D45').Value = '12.45'
$ws.Range('E45').Value = '  +3.60%  '
$ws.Range('D46').Value = '3.569'
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('D47').Value = '120.14'
$ws.Range('E47').Value = '  +8.87%  '
$ws.Range('D48').Value = '0.5344'
$ws.Range('E48').Value = '  +4.91%  '
$ws.Range('D49').Value = '1.808'
$ws.Range('E49').Value = '  +1.99%  '
$ws.Range('D50').Value = '1.053'
$ws.Range('E50').Value = '  +0.67%  '
$ws.Range('D51').Value = '0.06344'
$ws.Range('E51').Value = '  +2.04%  '
